# Insert a new data row at row 88 (pushing existing rows 88-193 down to 89-194)
# and populate it with the new weekly price record, per the authoring diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 88; this shifts rows 88-193
# down to 89-194 and extends the sheet dimension from A1:R193 to A1:R194.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new record's data.
$ws.Cells.Item(88, 1).Value = 9
$ws.Cells.Item(88, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(88, 3).Value = "Metropolitana"
$ws.Cells.Item(88, 4).Value = 44546
$ws.Cells.Item(88, 5).Value = 13
$ws.Cells.Item(88, 6).Value = 100112026
$ws.Cells.Item(88, 7).Value = "Haba"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 106
$ws.Cells.Item(88, 11).Value = 16000
$ws.Cells.Item(88, 12).Value = 17000
$ws.Cells.Item(88, 13).Value = 16500
$ws.Cells.Item(88, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(88, 15).Value = "Carahue"
$ws.Cells.Item(88, 16).Value = 660
$ws.Cells.Item(88, 17).Value = 25
$ws.Cells.Item(88, 18).Value = "Hortaliza"
